$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '27.103.21'
$ws.Cells.Item(2, 5).Value = '  +0.82%  '
$ws.Cells.Item(3, 4).Value = '1.812.70'
$ws.Cells.Item(3, 5).Value = '  +0.52%  '
$ws.Cells.Item(4, 5).Value = '  +0.16%  '
$ws.Cells.Item(5, 4).Value = '''311.48'
$ws.Cells.Item(5, 5).Value = '  +0.49%  '
$ws.Cells.Item(6, 5).Value = '  +0.17%  '
$ws.Cells.Item(7, 4).Value = '''0.4621'
$ws.Cells.Item(7, 5).Value = '  +4.85%  '
$ws.Cells.Item(8, 4).Value = '''0.3756'
$ws.Cells.Item(8, 5).Value = '  +1.84%  '
$ws.Cells.Item(9, 4).Value = '''0.07408'
$ws.Cells.Item(9, 5).Value = '  -0.11%  '
$ws.Cells.Item(10, 4).Value = '''0.8623'
$ws.Cells.Item(10, 5).Value = '  +0.51%  '
$ws.Cells.Item(11, 4).Value = '''20.59'
$ws.Cells.Item(11, 5).Value = '  -0.74%  '
$ws.Cells.Item(12, 4).Value = '1.815.56'
$ws.Cells.Item(12, 5).Value = '  +0.75%  '
$ws.Cells.Item(13, 4).Value = '''6.649'
$ws.Cells.Item(13, 5).Value = '  +0.36%  '
$ws.Cells.Item(14, 4).Value = '''5.385'
$ws.Cells.Item(14, 5).Value = '  +2.18%  '
$ws.Cells.Item(15, 2).Value = 'TRON'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Cells.Item(15, 4).Value = '''0.07076'
$ws.Cells.Item(15, 5).Value = '  +0.14%  '
$ws.Cells.Item(16, 2).Value = 'Litecoin'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(16, 4).Value = '''92.00'
$ws.Cells.Item(16, 5).Value = '  -1.05%  '
$ws.Cells.Item(17, 4).Value = '''1.001'
$ws.Cells.Item(17, 5).Value = '  +0.12%  '
$ws.Cells.Item(18, 4).Value = '''0.000008728'
$ws.Cells.Item(18, 5).Value = '  +0.44%  '
$ws.Cells.Item(19, 5).Value = '  +0.20%  '
$ws.Cells.Item(20, 4).Value = '''14.87'
$ws.Cells.Item(20, 5).Value = '  +0.44%  '
$ws.Cells.Item(21, 4).Value = '27.119.73'
$ws.Cells.Item(21, 5).Value = '  +0.80%  '
$ws.Cells.Item(22, 5).Value = '  +2.93%  '
$ws.Cells.Item(23, 5).Value = '  +0.35%  '
$ws.Cells.Item(24, 4).Value = '2.045.10'
$ws.Cells.Item(24, 5).Value = '  +1.17%  '
$ws.Cells.Item(25, 4).Value = '''1.923'
$ws.Cells.Item(25, 5).Value = '  -2.52%  '
$ws.Cells.Item(26, 4).Value = '''151.41'
$ws.Cells.Item(26, 5).Value = '  -0.07%  '
$ws.Cells.Item(27, 4).Value = '''2.200'
$ws.Cells.Item(27, 5).Value = '  -0.82%  '
$ws.Cells.Item(28, 4).Value = '''18.53'
$ws.Cells.Item(28, 5).Value = '  +0.86%  '
$ws.Cells.Item(29, 4).Value = '''5.264'
$ws.Cells.Item(29, 5).Value = '  +1.16%  '
$ws.Cells.Item(30, 4).Value = '''116.99'
$ws.Cells.Item(30, 5).Value = '  -0.39%  '
$ws.Cells.Item(31, 4).Value = '''0.08917'
$ws.Cells.Item(31, 5).Value = '  +1.53%  '
$ws.Cells.Item(32, 4).Value = '''0.7727'
$ws.Cells.Item(32, 5).Value = '  +4.06%  '
$ws.Cells.Item(33, 5).Value = '  +0.73%  '
$ws.Cells.Item(34, 5).Value = '  +1.09%  '
$ws.Cells.Item(35, 4).Value = '''2.899'
$ws.Cells.Item(35, 5).Value = '  +0.37%  '
$ws.Cells.Item(36, 4).Value = '''1.001'
$ws.Cells.Item(36, 5).Value = '  +0.18%  '
$ws.Cells.Item(37, 5).Value = '  +3.31%  '
$ws.Cells.Item(38, 5).Value = '  -0.40%  '
$ws.Cells.Item(39, 5).Value = '  +0.43%  '
$ws.Cells.Item(40, 2).Value = 'MXToken'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Cells.Item(40, 4).Value = '''2.929'
$ws.Cells.Item(40, 5).Value = '  +3.98%  '
$ws.Cells.Item(41, 2).Value = 'FraxShare'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Cells.Item(41, 4).Value = '''7.259'
$ws.Cells.Item(41, 5).Value = '  +2.69%  '
$ws.Cells.Item(42, 5).Value = '  +16.10%  '
$ws.Cells.Item(43, 4).Value = '''0.5276'
$ws.Cells.Item(43, 5).Value = '  +0.53%  '
$ws.Cells.Item(44, 4).Value = '''0.1678'
$ws.Cells.Item(44, 5).Value = '  -0.22%  '
$ws.Cells.Item(45, 4).Value = '''8.611'
$ws.Cells.Item(45, 5).Value = '  +1.51%  '
$ws.Cells.Item(46, 4).Value = '''0.5026'
$ws.Cells.Item(46, 5).Value = '  +0.53%  '
$ws.Cells.Item(47, 4).Value = '''10.32'
$ws.Cells.Item(47, 5).Value = '  -0.32%  '
$ws.Cells.Item(48, 4).Value = '''104.89'
$ws.Cells.Item(48, 5).Value = '  +0.54%  '
$ws.Cells.Item(49, 2).Value = 'PaxDollar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Cells.Item(49, 4).Value = '''1.001'
$ws.Cells.Item(49, 5).Value = '  +0.19%  '
$ws.Cells.Item(50, 2).Value = 'NEARProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(50, 4).Value = '''1.672'
$ws.Cells.Item(50, 5).Value = '  +0.24%  '
$ws.Cells.Item(51, 4).Value = '''0.06319'
$ws.Cells.Item(51, 5).Value = '  +0.14%  '
